$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update solenoid / motor channel assignment labels to match Robot Pinout
$ws.Range("A18").Value = "Ball Intake Motor"
$ws.Range("A34").Value = "Ball Arm Raise"
$ws.Range("A35").Value = "Ball Arm Lower"
$ws.Range("A36").Value = "Frame Stand Raise"
$ws.Range("A37").Value = "Frame Stand Lower"
